$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.070.93"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.876.17"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.45"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9984"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4919"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2922"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06616"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.886.50"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.56"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07207"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6673"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "86.23"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.914"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.037.67"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007816"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9985"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.83"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.124.95"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.789"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.859"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +4.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.157"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.57"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "143.03"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +9.77%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.898"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.386"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.217"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08793"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.997"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05089"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7227"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.113"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.657"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01859"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +12.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.682"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.166"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9286"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.784"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4243"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9978"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.12"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.386"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1282"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05694"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "32.89"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.319"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.61%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3776"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.343"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.75%  "
